$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing "Conto corrente eliminato" row (row 4) down to row 5,
# preserving its values + formatting (style indexes).
$ws.Range("A4:J4").Copy($ws.Range("A5:J5"))
$ws.Rows("5:5").RowHeight = 63.75

# Overwrite row 4 with the new "Modifica Conto corrente" test case.
$ws.Range("A4").Value = "Matrix Web : Conti Correnti_Verifica Modifica Conto corrente"
$ws.Range("B4").Value = "Ricerca di un cliente random >tab: Dettaglio Anagrafica > subtab: Conti correnti"
$ws.Range("C4").Value = "Verificare che il Conto corrente sia stato modificato correttamente nella tabella del tab Conti correnti "

# Update the active selection to match the edited row.
$ws.Range("A4:J4").Select()
